# emissions_cost.xlsx — "env values created + update of all tables according to new style"
#
# The meaningful content change: the "marginal" sheet's class column (A) was
# originally only populated on the first row of each class-group (merged look),
# rows 3/5/6/7/9/10/11 had no value in column A. The updated table fills the
# class label down onto every row so each row is fully self-describing.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("marginal")

# Fill the "class" column (A) down within each group so every row carries
# its class label (Short-haul / Medium-haul / Long-haul).
$ws.Range("A3").Value = $ws.Range("A2").Value2

$ws.Range("A5").Value = $ws.Range("A4").Value2
$ws.Range("A6").Value = $ws.Range("A4").Value2
$ws.Range("A7").Value = $ws.Range("A4").Value2

$ws.Range("A9").Value  = $ws.Range("A8").Value2
$ws.Range("A10").Value = $ws.Range("A8").Value2
$ws.Range("A11").Value = $ws.Range("A8").Value2

# Bring the "marginal" sheet to the front (it becomes the active/selected
# tab) and leave the cursor parked on F16, matching the saved view state.
$ws.Select() | Out-Null
$ws.Range("F16").Select() | Out-Null
